$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("D1").Value = "Assigned To"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Row 3 - new test case TC2 (ID + description first)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "TC2"
$ws.Range("C3").Value = "Website  should display Trending tags on search "

# Row 4 - new test case TC3 (ID only for now)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "TC3"

# Assigned To values
$ws.Range("D2").Value = "Rahul"
$ws.Range("D3").Value = "Satish"
$ws.Range("D4").Value = "Satish"

# Row 4 description (entered last)
$ws.Range("C4").Value = "When customer click  search button and select any trending tags then `n1) Flat offer pop up should be displayed`n2) customer is asked to specify if its for Girl or a boy`n3) customer must be given option of Tiny preemie, Preemie and 0-3M 3-5M,6-9M and +More`n4) Apply button is available`n5) When customer complete step 2-4 then min 5 products are displayed"
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(2).RowHeight

# Column widths (target OOXML width attrs: C=49, D=20)
$ws.Columns.Item(3).ColumnWidth = 48.166666666666664
$ws.Columns.Item(4).ColumnWidth = 19.166666666666668

# Selection
$ws.Range("C4").Select()
